$wb = $excel.ActiveWorkbook

# --- Resumen sheet: update Zona (max) label and Maximo value ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z4"
$wsResumen.Range("C2").Value = 540.4089361012002

# --- Solucion sheet: new randomized Pedido -> Salida assignment ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$solucionData = @(
    @("Pedido_27", "S001"),
    @("Pedido_59", "S021"),
    @("Pedido_26", "S041"),
    @("Pedido_43", "S061"),
    @("Pedido_34", "S011"),
    @("Pedido_30", "S071"),
    @("Pedido_77", "S051"),
    @("Pedido_11", "S031"),
    @("Pedido_48", "S062"),
    @("Pedido_20", "S042"),
    @("Pedido_61", "S002"),
    @("Pedido_69", "S022"),
    @("Pedido_2", "S072"),
    @("Pedido_36", "S032"),
    @("Pedido_5", "S012"),
    @("Pedido_64", "S052"),
    @("Pedido_35", "S063"),
    @("Pedido_63", "S023"),
    @("Pedido_66", "S043"),
    @("Pedido_33", "S003"),
    @("Pedido_74", "S033"),
    @("Pedido_14", "S073"),
    @("Pedido_21", "S053"),
    @("Pedido_65", "S064"),
    @("Pedido_72", "S013"),
    @("Pedido_12", "S024"),
    @("Pedido_42", "S044"),
    @("Pedido_40", "S004"),
    @("Pedido_29", "S074"),
    @("Pedido_56", "S054"),
    @("Pedido_54", "S014"),
    @("Pedido_32", "S034"),
    @("Pedido_28", "S065"),
    @("Pedido_8", "S005"),
    @("Pedido_44", "S045"),
    @("Pedido_38", "S025"),
    @("Pedido_41", "S055"),
    @("Pedido_55", "S075"),
    @("Pedido_18", "S015"),
    @("Pedido_3", "S035"),
    @("Pedido_23", "S046"),
    @("Pedido_10", "S066"),
    @("Pedido_57", "S006"),
    @("Pedido_1", "S016"),
    @("Pedido_75", "S056"),
    @("Pedido_19", "S026"),
    @("Pedido_49", "S007"),
    @("Pedido_47", "S076"),
    @("Pedido_80", "S047"),
    @("Pedido_37", "S036"),
    @("Pedido_58", "S067"),
    @("Pedido_73", "S017"),
    @("Pedido_70", "S057"),
    @("Pedido_78", "S077"),
    @("Pedido_45", "S008"),
    @("Pedido_62", "S027"),
    @("Pedido_50", "S048"),
    @("Pedido_79", "S068"),
    @("Pedido_51", "S037"),
    @("Pedido_25", "S058"),
    @("Pedido_52", "S018"),
    @("Pedido_39", "S078"),
    @("Pedido_68", "S028"),
    @("Pedido_24", "S049"),
    @("Pedido_53", "S038"),
    @("Pedido_76", "S009"),
    @("Pedido_60", "S069"),
    @("Pedido_4", "S059"),
    @("Pedido_16", "S029"),
    @("Pedido_13", "S039"),
    @("Pedido_46", "S079"),
    @("Pedido_67", "S019"),
    @("Pedido_9", "S050"),
    @("Pedido_71", "S010"),
    @("Pedido_15", "S070"),
    @("Pedido_17", "S030"),
    @("Pedido_7", "S060"),
    @("Pedido_31", "S040"),
    @("Pedido_22", "S020"),
    @("Pedido_6", "S080")
)

for ($i = 0; $i -lt $solucionData.Count; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucionData[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucionData[$i][1]
}

# --- Metricas sheet: updated Tiempo values per Zona ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 535.5234349659422
$wsMetricas.Range("B3").Value = 540.1137204022057
$wsMetricas.Range("B4").Value = 529.091426100119
$wsMetricas.Range("B5").Value = 540.4089361012002

Write-Host "Edit applied successfully"
